# Data refresh on the "Estado de Cuenta" sheet: previous contribution
# periods/amounts are removed and new ones are recorded for the
# NIT-9012445820 account.
#
# Rows 16 and 18 swap their "Periodo Mora" / "Valor Mora" pairs:
#   Row 16: 2107 / 36341  ->  2201 / 30284
#   Row 18: 2201 / 30284  ->  2107 / 36341
# Row 17 (2108 / 36341) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2201"
$ws.Range("F16").Value = 30284

$ws.Range("E18").Value = "2107"
$ws.Range("F18").Value = 36341
